$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether the value must be
# forced to stay text (numeric-looking strings that Excel would
# otherwise auto-convert to a number).
$updates = @(
    @{ Cell = 'D2'; Value = '63.669.80'; ForceText = 0 }
    @{ Cell = 'E2'; Value = '  +3.19%  '; ForceText = 0 }
    @{ Cell = 'D3'; Value = '3.128.93'; ForceText = 0 }
    @{ Cell = 'E3'; Value = '  +1.60%  '; ForceText = 0 }
    @{ Cell = 'E4'; Value = '  +0.00%  '; ForceText = 0 }
    @{ Cell = 'D5'; Value = '590.45'; ForceText = 1 }
    @{ Cell = 'E5'; Value = '  +2.24%  '; ForceText = 0 }
    @{ Cell = 'D6'; Value = '146.07'; ForceText = 1 }
    @{ Cell = 'E6'; Value = '  +2.96%  '; ForceText = 0 }
    @{ Cell = 'E7'; Value = '  +0.01%  '; ForceText = 0 }
    @{ Cell = 'D8'; Value = '3.120.19'; ForceText = 0 }
    @{ Cell = 'E8'; Value = '  +1.55%  '; ForceText = 0 }
    @{ Cell = 'E9'; Value = '  +1.84%  '; ForceText = 0 }
    @{ Cell = 'E10'; Value = '  +16.66%  '; ForceText = 0 }
    @{ Cell = 'E11'; Value = '  +3.96%  '; ForceText = 0 }
    @{ Cell = 'E12'; Value = '  +0.39%  '; ForceText = 0 }
    @{ Cell = 'D13'; Value = '0.0000254'; ForceText = 1 }
    @{ Cell = 'E13'; Value = '  +6.28%  '; ForceText = 0 }
    @{ Cell = 'D14'; Value = '36.13'; ForceText = 1 }
    @{ Cell = 'E14'; Value = '  +3.03%  '; ForceText = 0 }
    @{ Cell = 'E15'; Value = '  -0.45%  '; ForceText = 0 }
    @{ Cell = 'D16'; Value = '3.644.37'; ForceText = 0 }
    @{ Cell = 'E16'; Value = '  +1.58%  '; ForceText = 0 }
    @{ Cell = 'E17'; Value = '  -1.02%  '; ForceText = 0 }
    @{ Cell = 'D18'; Value = '63.596.04'; ForceText = 0 }
    @{ Cell = 'E18'; Value = '  +3.18%  '; ForceText = 0 }
    @{ Cell = 'D19'; Value = '3.131.31'; ForceText = 0 }
    @{ Cell = 'E19'; Value = '  +2.00%  '; ForceText = 0 }
    @{ Cell = 'D20'; Value = '465.97'; ForceText = 1 }
    @{ Cell = 'E20'; Value = '  +3.70%  '; ForceText = 0 }
    @{ Cell = 'D21'; Value = '14.24'; ForceText = 1 }
    @{ Cell = 'E21'; Value = '  +2.23%  '; ForceText = 0 }
    @{ Cell = 'D22'; Value = '0.736'; ForceText = 1 }
    @{ Cell = 'E22'; Value = '  +0.69%  '; ForceText = 0 }
    @{ Cell = 'D23'; Value = '7.54'; ForceText = 1 }
    @{ Cell = 'E23'; Value = '  +1.35%  '; ForceText = 0 }
    @{ Cell = 'E24'; Value = '  -2.03%  '; ForceText = 0 }
    @{ Cell = 'D25'; Value = '82.24'; ForceText = 1 }
    @{ Cell = 'E25'; Value = '  +0.21%  '; ForceText = 0 }
    @{ Cell = 'E26'; Value = '  -0.33%  '; ForceText = 0 }
    @{ Cell = 'D27'; Value = '8.73'; ForceText = 1 }
    @{ Cell = 'E27'; Value = '  +8.26%  '; ForceText = 0 }
    @{ Cell = 'E28'; Value = '  +2.74%  '; ForceText = 0 }
    @{ Cell = 'E29'; Value = '  -0.50%  '; ForceText = 0 }
    @{ Cell = 'E30'; Value = '  -0.04%  '; ForceText = 0 }
    @{ Cell = 'D31'; Value = '6.83'; ForceText = 1 }
    @{ Cell = 'E31'; Value = '  +1.99%  '; ForceText = 0 }
    @{ Cell = 'D32'; Value = '27.09'; ForceText = 1 }
    @{ Cell = 'E32'; Value = '  +1.86%  '; ForceText = 0 }
    @{ Cell = 'E33'; Value = '  +0.78%  '; ForceText = 0 }
    @{ Cell = 'E34'; Value = '  +8.50%  '; ForceText = 0 }
    @{ Cell = 'E35'; Value = '  +10.21%  '; ForceText = 0 }
    @{ Cell = 'E36'; Value = '  +1.66%  '; ForceText = 0 }
    @{ Cell = 'B37'; Value = 'Filecoin'; ForceText = 0 }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText = 0 }
    @{ Cell = 'D37'; Value = '6.13'; ForceText = 1 }
    @{ Cell = 'E37'; Value = '  +1.06%  '; ForceText = 0 }
    @{ Cell = 'B38'; Value = 'dogwifhat'; ForceText = 0 }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; ForceText = 0 }
    @{ Cell = 'D38'; Value = '3.35'; ForceText = 1 }
    @{ Cell = 'E38'; Value = '  +13.42%  '; ForceText = 0 }
    @{ Cell = 'D39'; Value = '50.87'; ForceText = 1 }
    @{ Cell = 'E39'; Value = '  +1.48%  '; ForceText = 0 }
    @{ Cell = 'D40'; Value = '447.40'; ForceText = 1 }
    @{ Cell = 'E40'; Value = '  +5.07%  '; ForceText = 0 }
    @{ Cell = 'E41'; Value = '  -0.96%  '; ForceText = 0 }
    @{ Cell = 'E42'; Value = '  +1.09%  '; ForceText = 0 }
    @{ Cell = 'D43'; Value = '2.920.32'; ForceText = 0 }
    @{ Cell = 'E43'; Value = '  +4.98%  '; ForceText = 0 }
    @{ Cell = 'D44'; Value = '0.278'; ForceText = 1 }
    @{ Cell = 'E44'; Value = '  +3.49%  '; ForceText = 0 }
    @{ Cell = 'E45'; Value = '  +2.33%  '; ForceText = 0 }
    @{ Cell = 'E46'; Value = '  +3.16%  '; ForceText = 0 }
    @{ Cell = 'D47'; Value = '125.64'; ForceText = 1 }
    @{ Cell = 'E47'; Value = '  +0.52%  '; ForceText = 0 }
    @{ Cell = 'B48'; Value = 'Arweave'; ForceText = 0 }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'; ForceText = 0 }
    @{ Cell = 'D48'; Value = '35.16'; ForceText = 1 }
    @{ Cell = 'E48'; Value = '  -0.54%  '; ForceText = 0 }
    @{ Cell = 'B49'; Value = 'USDe'; ForceText = 0 }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'; ForceText = 0 }
    @{ Cell = 'D49'; Value = '0.999'; ForceText = 1 }
    @{ Cell = 'E49'; Value = '  +0.00%  '; ForceText = 0 }
    @{ Cell = 'E50'; Value = '  +0.19%  '; ForceText = 0 }
    @{ Cell = 'D51'; Value = '24.63'; ForceText = 1 }
    @{ Cell = 'E51'; Value = '  +3.03%  '; ForceText = 0 }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText -eq 1) {
        # Numeric-looking text (e.g. "590.45") would otherwise be coerced
        # into a real number by Excel; force text, write it, then drop the
        # temporary number format so the cell keeps its original (default)
        # style, matching how the sheet stored these as inline strings.
        $range.NumberFormat = '@'
        $range.Value = $u.Value
        $range.ClearFormats()
    } else {
        $range.Value = $u.Value
    }
}
